$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$filesTabQuery = @'
MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)<-[*]-(prt)<--(f:file)
    WHERE c.gender = "FEMALE"
WITH DISTINCT f, prt, c, a, ct
RETURN
    COALESCE(f.file_name, '') AS `File Name`,
    COALESCE(head(labels(prt)), '') AS `Association`,
​
    COALESCE(f.file_description, '') AS `Description`,
    COALESCE(f.file_format, '') AS `File Format`,
    COALESCE(f.file_size, '') AS `Size`,
    COALESCE(ct.clinical_trial_designation, '') AS `Trial Code`,
    COALESCE(a.arm_id, '') AS `Arm`,
    COALESCE(c.case_id, '') AS `Case ID`
'@

$ws.Range("A3").Value = "FilesTab"
$ws.Range("B3").Value = $filesTabQuery
$ws.Range("B3").WrapText = $true
$ws.Range("C3").Value = $ws.Range("C2").Value2
$ws.Range("C3").WrapText = $true
$ws.Range("D3").Value = $ws.Range("D2").Value2
$ws.Range("E3").Value = $ws.Range("E2").Value2

$ws.Rows(3).RowHeight = 188.5

$ws.Range("D7").Select()
